$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge-field placeholder renames.
#    The legacy "object.field" dotted placeholders for the flattened
#    company / buyer / execution / delivery fields are renamed to a single
#    camelCase token (no dot). The "product.*" placeholders (which describe
#    a repeating line-item) are pluralised to "products.*" but keep their dot.
# ---------------------------------------------------------------------------

$renames = @(
    @("company.address", "companyAddress"),
    @("company.email",   "companyEmail"),
    @("company.phone",   "companyPhone"),
    @("company.name",    "companyName"),
    @("buyer.name",      "buyerName"),
    @("buyer.address",   "buyerAddress"),
    @("execution.date",  "executionDate"),
    @("delivery.date",   "deliveryDate"),
    @("product.name",    "products.name"),
    @("product.quantity","products.quantity"),
    @("product.price",   "products.price"),
    @("product.cost",    "products.cost")
)

foreach ($pair in $renames) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 2. Tidy up a couple of paragraphs whose sentences were artificially split
#    into extra runs around now-removed grammar-check markers. Re-typing the
#    same text over the existing range collapses them back into a single run.
# ---------------------------------------------------------------------------

$inspectionText = " Buyer is entitled to inspect the Goods upon delivery. If the Goods are unacceptable for any reason, Buyer must reject them at the time of delivery up to five (5) business days from the date of delivery. If Buyer has not rejected the Goods within five (5) business days from the date of delivery, Buyer shall have waived any right to reject that specific delivery of Goods. In the event Buyer rejects the Goods, Buyer shall allow Seller a reasonable time to cure the deficiency. A reasonable time period shall be determined by industry standards for the particular Goods, as well as the Seller and Buyer."
$d.Content.Find.Execute($inspectionText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $inspectionText, 2) | Out-Null

$governingLawText = "The Parties agree that this Agreement shall be governed by the State and/or Country in which both Parties do business. In the event that the Parties do business in different States and/or Countries, this Agreement shall be governed by _________ law."
$d.Content.Find.Execute($governingLawText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $governingLawText, 2) | Out-Null

$entireAgreementText = " The Parties acknowledge and agree that this Agreement represents the entire agreement between the Parties. In the event that the Parties desire to change, add, or otherwise modify any terms, they shall do so in writing to be signed by both parties."
$d.Content.Find.Execute($entireAgreementText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $entireAgreementText, 2) | Out-Null

# ---------------------------------------------------------------------------
# 2b. The "_GoBack" bookmark that used to sit inside "{{product.cost}}" lands
#     inside "{{deliveryDate}}" after the edits above (it tracks the last
#     place text was typed/replaced) - recreate it there.
# ---------------------------------------------------------------------------

$goBackRange = $d.Content
$found = $goBackRange.Find.Execute("deliveryDate", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if ($found) {
    $afterD = $goBackRange.Start + 9
    $bmRange = $d.Range($afterD, $afterD)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

# ---------------------------------------------------------------------------
# 3. Column width tweaks on the letterhead table and the line-items table.
# ---------------------------------------------------------------------------

$letterheadTable = $d.Tables.Item(1)
$letterheadTable.Columns.Item(1).Width = 366.2   # 7324 twips
$letterheadTable.Columns.Item(2).Width = 101.55  # 2031 twips

$lineItemsTable = $d.Tables.Item(3)
$lineItemsTable.Columns.Item(1).Width = 177.2    # 3544 twips
$lineItemsTable.Columns.Item(2).Width = 106.3    # 2126 twips
$lineItemsTable.Columns.Item(3).Width = 92.15    # 1843 twips
$lineItemsTable.Columns.Item(4).Width = 91.6     # 1832 twips
